$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''23.038.23'

$ws.Range("D3").Value = '''1.599.15'
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").Value = '''1.001'
$ws.Range("E5").Value = '  -0.25%  '

$ws.Range("D6").Value = '''302.19'
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").Value = '''0.3781'
$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").Value = '''0.3650'
$ws.Range("E8").Value = '  -0.62%  '

$ws.Range("D9").Value = '''50.66'
$ws.Range("E9").Value = '  +1.07%  '

$ws.Range("D10").Value = '''1.253'
$ws.Range("E10").Value = '  -2.28%  '

$ws.Range("B11").Value = 'BinanceUSD'
$ws.Range("C11").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D11").Value = '''1.002'
$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = '''0.08139'
$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").Value = '''22.33'
$ws.Range("E13").Value = '  -2.57%  '

$ws.Range("D14").Value = '''6.580'
$ws.Range("E14").Value = '  -1.58%  '

$ws.Range("D15").Value = '''7.368'
$ws.Range("E15").Value = '  -2.76%  '

$ws.Range("D16").Value = '''0.00001245'
$ws.Range("E16").Value = '  -2.13%  '

$ws.Range("D17").Value = '''1.598.36'
$ws.Range("E17").Value = '  -0.35%  '

$ws.Range("D18").Value = '''91.87'
$ws.Range("E18").Value = '  +0.35%  '

$ws.Range("D19").Value = '''0.06826'
$ws.Range("E19").Value = '  -0.04%  '

$ws.Range("D20").Value = '''18.16'
$ws.Range("E20").Value = '  -2.47%  '

$ws.Range("D21").Value = '''6.520'
$ws.Range("E21").Value = '  -1.97%  '

$ws.Range("E22").Value = '  -0.25%  '

$ws.Range("D23").Value = '''13.01'
$ws.Range("E23").Value = '  -1.21%  '

$ws.Range("D24").Value = '''23.045.66'
$ws.Range("E24").Value = '  -1.13%  '

$ws.Range("D25").Value = '''2.363'
$ws.Range("E25").Value = '  -1.41%  '

$ws.Range("D26").Value = '''2.775'
$ws.Range("E26").Value = '  -6.73%  '

$ws.Range("D27").Value = '''21.06'
$ws.Range("E27").Value = '  -0.87%  '

$ws.Range("D28").Value = '''148.91'
$ws.Range("E28").Value = '  -1.36%  '

$ws.Range("D29").Value = '''5.252'
$ws.Range("E29").Value = '  -1.36%  '

$ws.Range("D30").Value = '''134.29'
$ws.Range("E30").Value = '  +0.79%  '

$ws.Range("D31").Value = '''2.359'
$ws.Range("E31").Value = '  -4.64%  '

$ws.Range("D32").Value = '''6.820'
$ws.Range("E32").Value = '  -9.19%  '

$ws.Range("D33").Value = '''1.777.21'
$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("D34").Value = '''0.9605'
$ws.Range("E34").Value = '  -1.60%  '

$ws.Range("D35").Value = '''0.07570'
$ws.Range("E35").Value = '  -2.72%  '

$ws.Range("D36").Value = '''10.31'
$ws.Range("E36").Value = '  +0.80%  '

$ws.Range("D37").Value = '''6.224'
$ws.Range("E37").Value = '  -2.36%  '

$ws.Range("D38").Value = '''0.02707'
$ws.Range("E38").Value = '  -3.40%  '

$ws.Range("D39").Value = '''0.2517'
$ws.Range("E39").Value = '  -1.63%  '

$ws.Range("D40").Value = '''0.08819'
$ws.Range("E40").Value = '  -0.87%  '

$ws.Range("D41").Value = '''1.362'
$ws.Range("E41").Value = '  -2.54%  '

$ws.Range("D42").Value = '''0.7040'
$ws.Range("E42").Value = '  -2.40%  '

$ws.Range("D43").Value = '''12.31'
$ws.Range("E43").Value = '  -4.46%  '

$ws.Range("D44").Value = '''15.22'
$ws.Range("E44").Value = '  -4.87%  '

$ws.Range("D45").Value = '''0.6606'
$ws.Range("E45").Value = '  -0.58%  '

$ws.Range("D46").Value = '''0.9993'
$ws.Range("E46").Value = '  -0.30%  '

$ws.Range("D47").Value = '''2.286'
$ws.Range("E47").Value = '  -1.74%  '

$ws.Range("E48").Value = '  +0.26%  '

$ws.Range("D49").Value = '''131.70'
$ws.Range("E49").Value = '  -0.05%  '

$ws.Range("D50").Value = '''0.07934'
$ws.Range("E50").Value = '  -1.54%  '

$ws.Range("D51").Value = '''1.219'
$ws.Range("E51").Value = '  +3.20%  '
